# Swap the order of "dnasr281@gmail.com, System" -> "System, dnasr281@gmail.com"
# in column G ("Recorded By") for every row where the value matches exactly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$oldValue = "dnasr281@gmail.com, System"
$newValue = "System, dnasr281@gmail.com"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

$changed = 0
for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G ("Recorded By")
    if ($cell.Text -eq $oldValue) {
        $cell.Value = $newValue
        $changed = $changed + 1
    }
}

Write-Host "Updated $changed cells in column G"
